$wb = $excel.ActiveWorkbook

# Sheet ALC, row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 920.40814
$ws.Range("J17").Value = 939.7659
$ws.Range("L17").Value = 2819.2977
$ws.Range("N17").Value = -3155.2977

# Sheet ALC, row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 629.8
$ws.Range("I98").Value = 629.8
$ws.Range("K98").Value = 629.8
$ws.Range("M98").Value = 868.2

# Sheet ALC, row 105
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 54194.8
$ws.Range("J105").Value = 54194.8
$ws.Range("L105").Value = 54194.8
$ws.Range("N105").Value = -61182.8

# Sheet ALC, row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 629.8
$ws.Range("I122").Value = 629.8
$ws.Range("K122").Value = 1889.4
$ws.Range("M122").Value = 560.6000000000001

# Sheet ALC, row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -84800

# Sheet ALC, row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1786.6
$ws.Range("J129").Value = 3350
$ws.Range("L129").Value = 10050
$ws.Range("N129").Value = -20050

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2409.6538
$ws.Range("I137").Value = 2114.5789
$ws.Range("J137").Value = 3210.5715
$ws.Range("K137").Value = 6343.736699999999
$ws.Range("L137").Value = 9631.7145
$ws.Range("M137").Value = -3793.736699999999
$ws.Range("N137").Value = -14731.7145

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2666.6
$ws.Range("I32").Value = 2560.2554
$ws.Range("K32").Value = 2560.2554
$ws.Range("M32").Value = -2273.2554

# Sheet ARM, row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 966.6667
$ws.Range("I88").Value = 900
$ws.Range("K88").Value = 900
$ws.Range("M88").Value = -494

# Sheet ARM, row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 966.6667
$ws.Range("I91").Value = 900
$ws.Range("K91").Value = 900
$ws.Range("M91").Value = 504

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2166.5186
$ws.Range("I102").Value = 1865.2307
$ws.Range("K102").Value = 1865.2307
$ws.Range("M102").Value = -243.2307000000001

# Sheet BSM, row 43
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 21743562
$ws.Range("I86").Value = 3690.158
$ws.Range("J86").Value = 125007950
$ws.Range("K86").Value = 3690.158
$ws.Range("L86").Value = 125007950
$ws.Range("M86").Value = -2567.158
$ws.Range("N86").Value = -125010196

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 21743562
$ws.Range("I89").Value = 3690.158
$ws.Range("J89").Value = 125007950
$ws.Range("K89").Value = 18450.79
$ws.Range("L89").Value = 625039750
$ws.Range("M89").Value = -12834.79
$ws.Range("N89").Value = -625050982

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2704
$ws.Range("J16").Value = 4341.6665
$ws.Range("L16").Value = 4341.6665
$ws.Range("N16").Value = -4915.6665

# Sheet CRP, row 44
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

# Sheet CRP, row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7799.3335
$ws.Range("J86").Value = 6699.5
$ws.Range("L86").Value = 6699.5
$ws.Range("N86").Value = -8945.5

# Sheet CRP, row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 7799.3335
$ws.Range("J89").Value = 6699.5
$ws.Range("L89").Value = 33497.5
$ws.Range("N89").Value = -44729.5

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2704
$ws.Range("J113").Value = 4341.6665
$ws.Range("L113").Value = 4341.6665
$ws.Range("N113").Value = -8681.666499999999

# Sheet CUL, row 13
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 881.2308
$ws.Range("J13").Value = 1204.75
$ws.Range("L13").Value = 3614.25
$ws.Range("N13").Value = -3950.25

# Sheet CUL, row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 4257.364
$ws.Range("J23").Value = 3906
$ws.Range("L23").Value = 11718
$ws.Range("N23").Value = -12188

# Sheet CUL, row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 198.44444
$ws.Range("I33").Value = 54
$ws.Range("J33").Value = 254
$ws.Range("K33").Value = 324
$ws.Range("L33").Value = 1524
$ws.Range("M33").Value = -41
$ws.Range("N33").Value = -2090

# Sheet CUL, row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 1153.6666
$ws.Range("J38").Value = 2048.2
$ws.Range("L38").Value = 6144.599999999999
$ws.Range("N38").Value = -6838.599999999999

# Sheet CUL, row 106
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 14000
$ws.Range("J106").Value = 14000
$ws.Range("L106").Value = 42000
$ws.Range("N106").Value = -43892

# Sheet CUL, row 126
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 3833.3333
$ws.Range("I126").Value = 3500
$ws.Range("K126").Value = 10500
$ws.Range("M126").Value = -5560

# Sheet GSM, row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8337422
$ws.Range("I70").Value = 22225520
$ws.Range("K70").Value = 22225520
$ws.Range("M70").Value = -22225250

# Sheet GSM, row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8337422
$ws.Range("I73").Value = 22225520
$ws.Range("K73").Value = 22225520
$ws.Range("M73").Value = -22224584

# Sheet GSM, row 103
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8071.357
$ws.Range("J113").Value = 11111.111
$ws.Range("L113").Value = 11111.111
$ws.Range("N113").Value = -15451.111

# Sheet LTW, row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1451.4783
$ws.Range("I55").Value = 726
$ws.Range("K55").Value = 726
$ws.Range("M55").Value = -553

# Sheet LTW, row 106
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 19987.5
$ws.Range("J106").Value = 19987.5
$ws.Range("L106").Value = 19987.5
$ws.Range("N106").Value = -22511.5

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4244.3
$ws.Range("I122").Value = 3938.111
$ws.Range("K122").Value = 11814.333
$ws.Range("M122").Value = -9364.332999999999

# Sheet WVR, row 25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

# Sheet WVR, row 45
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 9999
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 9999
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 9999
$ws.Range("N45").Value = -10981
$ws.Range("M45").ClearContents()

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9198
$ws.Range("J62").Value = 10563.077
$ws.Range("L62").Value = 10563.077
$ws.Range("N62").Value = -11811.077

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 9198
$ws.Range("J65").Value = 10563.077
$ws.Range("L65").Value = 52815.38499999999
$ws.Range("N65").Value = -59055.38499999999

# Sheet WVR, row 104
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 21500
$ws.Range("J104").Value = 21500
$ws.Range("L104").Value = 21500
$ws.Range("N104").Value = -28488
